$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Solar Panels" between Ram Disk (row 20) and Sun_Sensor (row 21)
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = "Solar Panels"
$ws.Cells.Item(21, 3).Value = "First Upload of files"

# Insert a new row for "Launch Adapter" between IHU (row 13) and LBandDownConverter (row 14)
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = "Launch Adapter"
$ws.Cells.Item(14, 3).Value = "Mechanical Drawings uploaded"

# Update Antennas status (row 4, column C) now that the antenna switch work is done
$ws.Range("C4").Value = "Ant Switch complete, Plots of ants added"

# Move the selection to C5, matching the author's edit location
$ws.Range("C5").Select()
